$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the numeric-looking receipt numbers in column B to Text format so they are
# stored as text, matching the source data (inline strings), not auto-converted to
# numbers. Only touch the specific cells that hold a receipt number - leave the
# other (blank) B cells completely untouched.
$ws.Range("B821:B823").NumberFormat = "@"
$ws.Range("B830:B832").NumberFormat = "@"

# Row 821
$ws.Cells.Item(821, 1).Value = 44320
$ws.Cells.Item(821, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(821, 2).Value = "5040739"
$ws.Cells.Item(821, 3).Value = 3011
$ws.Cells.Item(821, 4).Value = "Order 5040739 Card(Stripe)"
$ws.Cells.Item(821, 6).Value = 973.21

# Row 822
$ws.Cells.Item(822, 1).Value = 44320
$ws.Cells.Item(822, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(822, 2).Value = "5040739"
$ws.Cells.Item(822, 3).Value = 2611
$ws.Cells.Item(822, 4).Value = "Order 5040739 Card(Stripe)"
$ws.Cells.Item(822, 6).Value = 116.79

# Row 823
$ws.Cells.Item(823, 1).Value = 44320
$ws.Cells.Item(823, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(823, 2).Value = "5040739"
$ws.Cells.Item(823, 3).Value = 1930
$ws.Cells.Item(823, 4).Value = "Order 5040739 Card(Stripe)"
$ws.Cells.Item(823, 5).Value = 1090

# Row 824
$ws.Cells.Item(824, 1).Value = 44321
$ws.Cells.Item(824, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(824, 3).Value = 4010
$ws.Cells.Item(824, 4).Value = "TINGSTAD PAPPER"
$ws.Cells.Item(824, 5).Value = 3102.4

# Row 825
$ws.Cells.Item(825, 1).Value = 44321
$ws.Cells.Item(825, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(825, 3).Value = 2641
$ws.Cells.Item(825, 4).Value = "TINGSTAD PAPPER"
$ws.Cells.Item(825, 5).Value = 775.6

# Row 826
$ws.Cells.Item(826, 1).Value = 44321
$ws.Cells.Item(826, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(826, 3).Value = 1930
$ws.Cells.Item(826, 4).Value = "TINGSTAD PAPPER"
$ws.Cells.Item(826, 6).Value = 3878

# Row 827
$ws.Cells.Item(827, 1).Value = 44321
$ws.Cells.Item(827, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(827, 3).Value = 4010
$ws.Cells.Item(827, 4).Value = "NGROCERIES K0135"
$ws.Cells.Item(827, 5).Value = 213.39

# Row 828
$ws.Cells.Item(828, 1).Value = 44321
$ws.Cells.Item(828, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(828, 3).Value = 2645
$ws.Cells.Item(828, 4).Value = "NGROCERIES K0135"
$ws.Cells.Item(828, 5).Value = 25.61

# Row 829
$ws.Cells.Item(829, 1).Value = 44321
$ws.Cells.Item(829, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(829, 3).Value = 1930
$ws.Cells.Item(829, 4).Value = "NGROCERIES K0135"
$ws.Cells.Item(829, 6).Value = 239

# Row 830
$ws.Cells.Item(830, 1).Value = 44322
$ws.Cells.Item(830, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(830, 2).Value = "5061308"
$ws.Cells.Item(830, 3).Value = 3011
$ws.Cells.Item(830, 4).Value = "Order 5061308 Swish +46723656673"
$ws.Cells.Item(830, 6).Value = 806.25

# Row 831
$ws.Cells.Item(831, 1).Value = 44322
$ws.Cells.Item(831, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(831, 2).Value = "5061308"
$ws.Cells.Item(831, 3).Value = 2611
$ws.Cells.Item(831, 4).Value = "Order 5061308 Swish +46723656673"
$ws.Cells.Item(831, 6).Value = 96.75

# Row 832
$ws.Cells.Item(832, 1).Value = 44322
$ws.Cells.Item(832, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(832, 2).Value = "5061308"
$ws.Cells.Item(832, 3).Value = 1930
$ws.Cells.Item(832, 4).Value = "Order 5061308 Swish +46723656673"
$ws.Cells.Item(832, 5).Value = 903

# Row 833
$ws.Cells.Item(833, 1).Value = 44322
$ws.Cells.Item(833, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(833, 3).Value = 5010
$ws.Cells.Item(833, 4).Value = "April hyra"
$ws.Cells.Item(833, 5).Value = 4166

# Row 834
$ws.Cells.Item(834, 1).Value = 44322
$ws.Cells.Item(834, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(834, 4).Value = "April hyra"
$ws.Cells.Item(834, 5).Value = 0

# Row 835
$ws.Cells.Item(835, 1).Value = 44322
$ws.Cells.Item(835, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(835, 3).Value = 1930
$ws.Cells.Item(835, 4).Value = "April hyra"
$ws.Cells.Item(835, 6).Value = 4166

# Row 836
$ws.Cells.Item(836, 1).Value = 44325
$ws.Cells.Item(836, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(836, 3).Value = 5670
$ws.Cells.Item(836, 4).Value = "ST1 V#LLINGBY K0135"
$ws.Cells.Item(836, 5).Value = 844.58

# Row 837
$ws.Cells.Item(837, 1).Value = 44325
$ws.Cells.Item(837, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(837, 3).Value = 2641
$ws.Cells.Item(837, 4).Value = "ST1 V#LLINGBY K0135"
$ws.Cells.Item(837, 5).Value = 211.14

# Row 838
$ws.Cells.Item(838, 1).Value = 44325
$ws.Cells.Item(838, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(838, 3).Value = 1930
$ws.Cells.Item(838, 4).Value = "ST1 V#LLINGBY K0135"
$ws.Cells.Item(838, 6).Value = 1055.72
